$d = $word.ActiveDocument

# The "Achieved 87% ...", "Built redistricting platform ...", "Discovered
# systematic race coding errors ..." and "Trigonometric algorithm ..." lines
# are near-duplicated elsewhere in the resume (Professional Experience /
# Key Projects), so matching on leading text alone is ambiguous. Scope the
# edit to the "KEY ACHIEVEMENTS AND IMPACT" section specifically: find that
# Heading2 paragraph, then the next Heading2 (or end of document) bounds the
# section.

$sectionStart = -1
$sectionEnd = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $styleName = $p.Style.NameLocal

    if ($styleName -eq "Heading 2") {
        if ($sectionStart -ne -1 -and $sectionEnd -eq -1) {
            $sectionEnd = $i - 1
        }
        if ($sectionStart -eq -1 -and $p.Range.Text -like "*KEY ACHIEVEMENTS AND IMPACT*") {
            $sectionStart = $i
        }
    }
}
if ($sectionStart -ne -1 -and $sectionEnd -eq -1) {
    $sectionEnd = $d.Paragraphs.Count
}

# Within that section, locate the six bullet paragraphs by their distinctive
# leading text.
$idxAchieved   = -1
$idxDelivered  = -1
$idxBuilt      = -1
$idxDeveloped  = -1
$idxDiscovered = -1
$idxTrig       = -1

for ($i = $sectionStart; $i -le $sectionEnd; $i++) {
    $t = $d.Paragraphs($i).Range.Text

    if ($t -like "*Achieved 87% prediction accuracy*") {
        $idxAchieved = $i
    }
    elseif ($t -like "*Delivered*4.9M additional revenue through continuous testing*") {
        $idxDelivered = $i
    }
    elseif ($t -like "*Built redistricting platform*") {
        $idxBuilt = $i
    }
    elseif ($t -like "*Developed longitudinal data analysis methods*") {
        $idxDeveloped = $i
    }
    elseif ($t -like "*Discovered systematic race coding errors*") {
        $idxDiscovered = $i
    }
    elseif ($t -like "*Trigonometric algorithm for boundary estimation*") {
        $idxTrig = $i
    }
}

# Rewrite the four retained bullets in place.
$d.Paragraphs($idxAchieved).Range.Text  = "• Revenue generation: Delivered `$4.9M additional revenue through optimization"
$d.Paragraphs($idxDelivered).Range.Text = "• 23% conversion rate improvement"
$d.Paragraphs($idxBuilt).Range.Text     = "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"
$d.Paragraphs($idxDeveloped).Range.Text = "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"

# Remove the two trailing bullets entirely (delete the higher index first so
# the lower index remains valid).
$toDelete = @($idxDiscovered, $idxTrig) | Sort-Object -Descending
foreach ($idx in $toDelete) {
    $d.Paragraphs($idx).Range.Delete()
}
